# Convert legacy <w:fldSimple w:instr="..."/> fields into the modern
# begin/instrText/separate/end run-based field representation, exactly like
# Word rewrites them as soon as a user touches/edits the field in the UI
# (see commit: "the usercontent m2doc is generated with wrong word xml
# syntax" - editing with MS Word turns fldSimple into fldChar run pairs).
#
# Targets the two M2Doc user-content markers used by this template:
#   m:usercontent zone1
#   m:endusercontent
#
# Any markup that was already inside the hosting paragraph before the field
# (e.g. the hidden "_GoBack" bookmark on the "endusercontent" paragraph) is
# preserved, as are the paragraph's own rsid-ish attributes.

$d = $word.ActiveDocument

# Field instruction codes (trimmed) that must be converted, in document order.
$targets = @("m:usercontent zone1", "m:endusercontent")

function Get-DocumentXmlChunks($doc) {
    $full = $doc.WordOpenXML
    $docStart = $full.IndexOf("<w:document")
    $endMarker = "</w:document>"
    $docEnd = $full.IndexOf($endMarker, $docStart) + $endMarker.Length
    $docXml = $full.Substring($docStart, $docEnd - $docStart)
    return ($docXml -split '(?=<w:p[ >])')
}

function Get-FieldParagraphParts($paraChunks, $instr) {
    $needle = '<w:fldSimple w:instr="' + $instr + '"'
    foreach ($chunk in $paraChunks) {
        if ($chunk.Contains($needle)) {
            $attrMatch = [regex]::Match($chunk, '^<w:p([^>]*)>')
            $attrsRaw = $attrMatch.Groups[1].Value

            # Only keep genuine w:rsid* attributes (paraId/textId are
            # synthesised by WordOpenXML and are not part of the saved file).
            $rsidAttrs = [regex]::Matches($attrsRaw, 'w:rsid\w*="[0-9A-Fa-f]*"')
            $attrParts = @()
            foreach ($a in $rsidAttrs) { $attrParts += $a.Value }
            $attrSuffix = ""
            if ($attrParts.Count -gt 0) { $attrSuffix = " " + ($attrParts -join " ") }

            # Keep whatever markup precedes the field inside the paragraph
            # (e.g. bookmarkStart/bookmarkEnd) so it is not lost.
            $openEnd = $attrMatch.Index + $attrMatch.Length
            $fldIdx = $chunk.IndexOf($needle)
            $prefix = $chunk.Substring($openEnd, $fldIdx - $openEnd)

            return @{ AttrSuffix = $attrSuffix; Prefix = $prefix }
        }
    }
    return @{ AttrSuffix = ""; Prefix = "" }
}

foreach ($instr in $targets) {
    $paraChunks = Get-DocumentXmlChunks $d
    $parts = Get-FieldParagraphParts $paraChunks $instr

    # Re-resolve the field each time: earlier conversions change the
    # Fields collection / range offsets, so always work off fresh data.
    $target = $null
    foreach ($f in $d.Fields) {
        if ($f.Code.Text.Trim() -eq $instr) {
            $target = $f
            break
        }
    }

    if ($target -eq $null) {
        Write-Output ("Field not found for instr: " + $instr)
        continue
    }

    $para = $target.Code.Paragraphs(1)
    $range = $para.Range

    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"' + $parts.AttrSuffix + '>' +
           $parts.Prefix +
           '<w:r><w:fldChar w:fldCharType="begin"/></w:r>' +
           '<w:r><w:instrText>' + $instr + '</w:instrText></w:r>' +
           '<w:r><w:fldChar w:fldCharType="separate"/></w:r>' +
           '<w:r><w:fldChar w:fldCharType="end"/></w:r>' +
           '</w:p>'

    $range.InsertXML($xml)
    Write-Output ("Converted field: " + $instr)
}
